$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 data (numeric-looking values must stay as text, like the rest of the sheet)
$ws.Range("A4:H4").NumberFormat = "@"

$ws.Range("A4").Value = "5"
$ws.Range("B4").Value = "4"
$ws.Range("C4").Value = "60"
$ws.Range("D4").Value = "What level of income from a"
$ws.Range("E4").Value = "What level of income from a"
$ws.Range("F4").Value = "0"
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
